$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 203.5
$ws.Range("I11").Value = 203.5
$ws.Range("K11").Value = 203.5
$ws.Range("M11").Value = -63.5
$ws.Range("H28").Value = 79174.69500000001
$ws.Range("I28").Value = 102256.5
$ws.Range("J28").Value = 2235.3333
$ws.Range("K28").Value = 102256.5
$ws.Range("L28").Value = 2235.3333
$ws.Range("M28").Value = -101771.5
$ws.Range("N28").Value = -3205.3333
$ws.Range("H43").Value = 8040.1
$ws.Range("I43").Value = 16000
$ws.Range("J43").Value = 6050.125
$ws.Range("K43").Value = 16000
$ws.Range("L43").Value = 6050.125
$ws.Range("M43").Value = -15931
$ws.Range("N43").Value = -6188.125
$ws.Range("H69").Value = 7844.615
$ws.Range("I69").Value = 4799
$ws.Range("J69").Value = 9198.223
$ws.Range("K69").Value = 14397
$ws.Range("L69").Value = 27594.669
$ws.Range("M69").Value = -13523
$ws.Range("N69").Value = -29342.669
$ws.Range("H72").Value = 7844.615
$ws.Range("I72").Value = 4799
$ws.Range("J72").Value = 9198.223
$ws.Range("K72").Value = 43191
$ws.Range("L72").Value = 82784.007
$ws.Range("M72").Value = -38823
$ws.Range("N72").Value = -91520.007
$ws.Range("H116").Value = 86114.5
$ws.Range("I116").Value = 112486
$ws.Range("K116").Value = 112486
$ws.Range("M116").Value = -109044
$ws.Range("H125").Value = 5876.067
$ws.Range("I125").Value = 5980.6665
$ws.Range("J125").Value = 5806.3335
$ws.Range("K125").Value = 53825.9985
$ws.Range("L125").Value = 52257.0015
$ws.Range("M125").Value = -51365.9985
$ws.Range("N125").Value = -57177.0015
$ws.Range("H132").Value = 18787
$ws.Range("I132").Value = 18787
$ws.Range("K132").Value = 56361
$ws.Range("M132").Value = -53831
$ws.Range("H138").Value = 4480
$ws.Range("I138").Value = 2966.6667
$ws.Range("K138").Value = 8900.000100000001
$ws.Range("M138").Value = -3760.000100000001
$ws.Range("H141").Value = 1568.8182
$ws.Range("I141").Value = 1555.7
$ws.Range("J141").Value = 1700
$ws.Range("K141").Value = 4667.1
$ws.Range("L141").Value = 5100
$ws.Range("M141").Value = 512.8999999999996
$ws.Range("N141").Value = -15460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 700.90125
$ws.Range("I32").Value = 634.6625
$ws.Range("K32").Value = 634.6625
$ws.Range("M32").Value = -347.6625
$ws.Range("H53").Value = 23749.75
$ws.Range("I53").Value = 17500
$ws.Range("J53").Value = 29999.5
$ws.Range("K53").Value = 17500
$ws.Range("L53").Value = 29999.5
$ws.Range("M53").Value = -16818
$ws.Range("N53").Value = -31363.5
$ws.Range("H74").Value = 4469.9443
$ws.Range("I74").Value = 1800.9
$ws.Range("J74").Value = 7806.25
$ws.Range("K74").Value = 1800.9
$ws.Range("L74").Value = 7806.25
$ws.Range("M74").Value = -926.9000000000001
$ws.Range("N74").Value = -9554.25
$ws.Range("H77").Value = 4469.9443
$ws.Range("I77").Value = 1800.9
$ws.Range("J77").Value = 7806.25
$ws.Range("K77").Value = 9004.5
$ws.Range("L77").Value = 39031.25
$ws.Range("M77").Value = -4636.5
$ws.Range("N77").Value = -47767.25
$ws.Range("H102").Value = 111330.25
$ws.Range("I102").Value = 146807.33
$ws.Range("K102").Value = 146807.33
$ws.Range("M102").Value = -145185.33
$ws.Range("H132").Value = 2110.1765
$ws.Range("I132").Value = 2064.4167
$ws.Range("K132").Value = 6193.250100000001
$ws.Range("M132").Value = -3663.250100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2206.2727
$ws.Range("I80").Value = 547
$ws.Range("K80").Value = 547
$ws.Range("M80").Value = 451
$ws.Range("H83").Value = 2206.2727
$ws.Range("I83").Value = 547
$ws.Range("K83").Value = 2735
$ws.Range("M83").Value = 2257
$ws.Range("H107").Value = 1543.1555
$ws.Range("I107").Value = 1522.3448
$ws.Range("K107").Value = 1522.3448
$ws.Range("M107").Value = 397.6551999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12790.359
$ws.Range("I31").Value = 1153.52
$ws.Range("K31").Value = 1153.52
$ws.Range("M31").Value = -858.52
$ws.Range("H34").Value = 12790.359
$ws.Range("I34").Value = 1153.52
$ws.Range("K34").Value = 1153.52
$ws.Range("M34").Value = -951.52
$ws.Range("H107").Value = 628.86664
$ws.Range("I107").Value = 554.9
$ws.Range("J107").Value = 776.8
$ws.Range("K107").Value = 554.9
$ws.Range("L107").Value = 776.8
$ws.Range("M107").Value = 1365.1
$ws.Range("N107").Value = -4616.8
$ws.Range("H132").Value = 4527.5
$ws.Range("J132").Value = 12000
$ws.Range("L132").Value = 36000
$ws.Range("N132").Value = -41060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 76923416
$ws.Range("J17").Value = 250000320
$ws.Range("L17").Value = 750000960
$ws.Range("N17").Value = -750001298
$ws.Range("H97").Value = 288.625
$ws.Range("J97").Value = 267.8
$ws.Range("L97").Value = 803.4000000000001
$ws.Range("N97").Value = -1795.4
$ws.Range("H108").Value = 10851.637
$ws.Range("I108").Value = 4449.6
$ws.Range("K108").Value = 13348.8
$ws.Range("M108").Value = -10468.8
$ws.Range("H129").Value = 1371.4231
$ws.Range("I129").Value = 370.46667
$ws.Range("J129").Value = 2736.3635
$ws.Range("K129").Value = 1111.40001
$ws.Range("L129").Value = 8209.0905
$ws.Range("M129").Value = 3888.59999
$ws.Range("N129").Value = -18209.0905
$ws.Range("H131").Value = 3261
$ws.Range("J131").Value = 3460.25
$ws.Range("L131").Value = 10380.75
$ws.Range("N131").Value = -20460.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16327.866
$ws.Range("I70").Value = 7974.75
$ws.Range("J70").Value = 19365.363
$ws.Range("K70").Value = 7974.75
$ws.Range("L70").Value = 19365.363
$ws.Range("M70").Value = -7704.75
$ws.Range("N70").Value = -19905.363
$ws.Range("H73").Value = 16327.866
$ws.Range("I73").Value = 7974.75
$ws.Range("J73").Value = 19365.363
$ws.Range("K73").Value = 7974.75
$ws.Range("L73").Value = 19365.363
$ws.Range("M73").Value = -7038.75
$ws.Range("N73").Value = -21237.363
$ws.Range("H102").Value = 3272.88
$ws.Range("J102").Value = 4430.364
$ws.Range("L102").Value = 4430.364
$ws.Range("N102").Value = -7674.364
$ws.Range("H113").Value = 3666.5
$ws.Range("I113").Value = 2305.1667
$ws.Range("K113").Value = 2305.1667
$ws.Range("M113").Value = -135.1667000000002
$ws.Range("H122").Value = 3531.6
$ws.Range("I122").Value = 4592.5386
$ws.Range("J122").Value = 2382.25
$ws.Range("K122").Value = 13777.6158
$ws.Range("L122").Value = 7146.75
$ws.Range("M122").Value = -11327.6158
$ws.Range("N122").Value = -12046.75
$ws.Range("H123").Value = 42000
$ws.Range("J123").Value = 42000
$ws.Range("L123").Value = 42000
$ws.Range("N123").Value = -46900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8811.129000000001
$ws.Range("I22").Value = 14815.883
$ws.Range("J22").Value = 1519.6428
$ws.Range("K22").Value = 14815.883
$ws.Range("L22").Value = 1519.6428
$ws.Range("M22").Value = -14520.883
$ws.Range("N22").Value = -2109.6428
$ws.Range("H27").Value = 8811.129000000001
$ws.Range("I27").Value = 14815.883
$ws.Range("J27").Value = 1519.6428
$ws.Range("K27").Value = 14815.883
$ws.Range("L27").Value = 1519.6428
$ws.Range("M27").Value = -14708.883
$ws.Range("N27").Value = -1733.6428
$ws.Range("H47").Value = 35000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 35000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 35000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -35980
$ws.Range("H52").Value = 35000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 35000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 35000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -35466
$ws.Range("H68").Value = 2955
$ws.Range("I68").Value = 2799.2856
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 2799.2856
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -2050.2856
$ws.Range("N68").Value = -4998
$ws.Range("H71").Value = 2955
$ws.Range("I71").Value = 2799.2856
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 13996.428
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -10252.428
$ws.Range("N71").Value = -24988
$ws.Range("H95").Value = 67161.71000000001
$ws.Range("J95").Value = 67161.71000000001
$ws.Range("L95").Value = 67161.71000000001
$ws.Range("N95").Value = -72653.71000000001
$ws.Range("H122").Value = 68670.75
$ws.Range("I122").Value = 5773.6665
$ws.Range("J122").Value = 257362
$ws.Range("K122").Value = 17320.9995
$ws.Range("L122").Value = 772086
$ws.Range("M122").Value = -14870.9995
$ws.Range("N122").Value = -776986

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 60285.832
$ws.Range("J95").Value = 60285.832
$ws.Range("L95").Value = 60285.832
$ws.Range("N95").Value = -65777.83199999999
$ws.Range("H98").Value = 95446.71000000001
$ws.Range("J98").Value = 95446.71000000001
$ws.Range("L98").Value = 95446.71000000001
$ws.Range("N98").Value = -101436.71
$ws.Range("H107").Value = 3051.3333
$ws.Range("I107").Value = 2096.875
$ws.Range("J107").Value = 4142.143
$ws.Range("K107").Value = 6290.625
$ws.Range("L107").Value = 12426.429
$ws.Range("M107").Value = -4370.625
$ws.Range("N107").Value = -16266.429
$ws.Range("H115").Value = 98900
$ws.Range("J115").Value = 98900
$ws.Range("L115").Value = 98900
$ws.Range("N115").Value = -102034
$ws.Range("H132").Value = 3691.0908
$ws.Range("I132").Value = 2255.1892
$ws.Range("K132").Value = 6765.567599999999
$ws.Range("M132").Value = -4235.567599999999
$ws.Range("H136").Value = 14381.36
$ws.Range("I136").Value = 15758.421
$ws.Range("K136").Value = 47275.263
$ws.Range("M136").Value = -44725.263
